$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as TEXT (inline/shared string), even when
# the value looks like a number, without leaving a residual "Text" number
# format / style index behind on the cell (the source file has no per-cell
# styles on these data cells).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    # Leading apostrophe forces Excel to store numeric-looking text as text
    # rather than silently coercing it to a Double.
    $r.Value2 = "'" + $val
    # Excel auto-applies a quoted-text number format when it detects this;
    # reset the style back to Normal so the cell keeps its original
    # (unstyled) appearance while remaining text internally.
    $r.Style = "Normal"
}

Set-TextValue "D2" "244.85"
Set-TextValue "D3" "23.06"
Set-TextValue "D4" "5.416"
Set-TextValue "D5" "0.06032"
Set-TextValue "D6" "3.394"
Set-TextValue "D8" "0.9280"
Set-TextValue "D9" "0.1425"
Set-TextValue "D10" "0.07435"
Set-TextValue "D11" "0.03391"
Set-TextValue "D12" "0.03074"
Set-TextValue "D13" "0.09367"
Set-TextValue "D14" "3.944"
Set-TextValue "D15" "0.001590"
Set-TextValue "D16" "0.04835"
Set-TextValue "D17" "0.0005944"
Set-TextValue "D18" "0.005311"
Set-TextValue "D19" "0.004155"
Set-TextValue "D20" "0.0009846"
Set-TextValue "D21" "0.00008706"
Set-TextValue "D23" "6.441"
Set-TextValue "D25" "0.3244"
Set-TextValue "D27" "0.0002448"
Set-TextValue "D40" "0.03972"

# Rows 41-43 got reshuffled (Kick -> BKEX -> CEJI -> Kick) with refreshed data
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1075"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002712"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003033"
$ws.Range("E43").Value = "42KickTokenKICK"

Set-TextValue "D44" "0.005964"
Set-TextValue "D45" "0.00005202"
Set-TextValue "D47" "0.0005804"
Set-TextValue "D48" "0.8506"
